$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking values
# (e.g. "1.00", "155.10", "63.719.49") keep their exact text
# representation instead of being auto-converted to numbers.
$colD = $ws.Range("D2:D51")
$colD.NumberFormat = "@"

$ws.Range("D2").Value = "63.719.49"
$ws.Range("E2").Value = "  -1.21%  "
$ws.Range("D3").Value = "2.630.03"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "576.97"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").Value = "155.10"
$ws.Range("E6").Value = "  -0.62%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "0.623"
$ws.Range("E8").Value = "  -3.42%  "
$ws.Range("D9").Value = "2.627.61"
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("E10").Value = "  -3.86%  "
$ws.Range("D11").Value = "5.81"
$ws.Range("E11").Value = "  +0.06%  "
$ws.Range("E12").Value = "  -2.09%  "
$ws.Range("E13").Value = "  +0.80%  "
$ws.Range("E14").Value = "  -0.70%  "
$ws.Range("D15").Value = "3.106.01"
$ws.Range("E15").Value = "  +0.28%  "
$ws.Range("E16").Value = "  -2.35%  "
$ws.Range("D17").Value = "63.666.34"
$ws.Range("E17").Value = "  -1.09%  "
$ws.Range("D18").Value = "2.644.81"
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("E19").Value = "  -1.21%  "
$ws.Range("D20").Value = "7.63"
$ws.Range("E20").Value = "  +3.42%  "
$ws.Range("E21").Value = "  -3.44%  "
$ws.Range("D22").Value = "344.55"
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("E23").Value = "  +0.37%  "
$ws.Range("D24").Value = "67.94"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("E25").Value = "  +9.03%  "
$ws.Range("E26").Value = "  -4.49%  "
$ws.Range("D27").Value = "604.96"
$ws.Range("E27").Value = "  +9.21%  "
$ws.Range("D28").Value = "9.25"
$ws.Range("E28").Value = "  -1.32%  "
$ws.Range("E29").Value = "  +2.57%  "
$ws.Range("D30").Value = "7.96"
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("E31").Value = "  -0.48%  "
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("E33").Value = "  -0.49%  "
$ws.Range("D34").Value = "1.76"
$ws.Range("E34").Value = "  +1.11%  "
$ws.Range("D35").Value = "6.64"
$ws.Range("E35").Value = "  +3.23%  "
$ws.Range("E36").Value = "  +0.88%  "
$ws.Range("E37").Value = "  -2.52%  "
$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").Value = "19.72"
$ws.Range("E39").Value = "  -1.55%  "
$ws.Range("E40").Value = "  -2.19%  "
$ws.Range("D41").Value = "150.07"
$ws.Range("E41").Value = "  -1.10%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("E43").Value = "  +2.51%  "
$ws.Range("D44").Value = "41.74"
$ws.Range("E44").Value = "  -0.66%  "
$ws.Range("D45").Value = "24.79"
$ws.Range("E45").Value = "  +8.37%  "
$ws.Range("D46").Value = "159.21"
$ws.Range("E46").Value = "  +0.54%  "
$ws.Range("E47").Value = "  -2.12%  "
$ws.Range("E48").Value = "  -2.46%  "
$ws.Range("E49").Value = "  -0.54%  "
$ws.Range("D50").Value = "0.0997"
$ws.Range("E50").Value = "  -1.75%  "
$ws.Range("E51").Value = "  -0.70%  "

# Restore the default (unstyled) cell style on column D now that the
# text values are safely stored, matching the original formatting.
$colD.Style = "Normal"
